$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 32 for the new "strWindowPos" setting entry.
$ws.Rows.Item(32).Insert()

# Column B (File) keeps the same shared value as the rest of the table rows.
$ws.Cells.Item(32, 2).Value = $ws.Cells.Item(33, 2).Value()

# New shared strings must be created in this exact order so they land on the
# same shared-string table indices as in the target workbook (391, 392, 393):
#   391 = strWindowPos
#   392 = In "settings" form, tab "User interface"
#   393 = Remember window position and size on startup
$ws.Cells.Item(32, 3).Value = "strWindowPos"
$ws.Cells.Item(32, 4).Value = "In ""settings"" form, tab ""User interface"""
$ws.Cells.Item(32, 5).Value = "Remember window position and size on startup"

# Set the Comment (column D) for the existing "strChkDlgPath" row (now row 25,
# unaffected by the insert above) to reference the same "settings" form /
# "User interface" tab comment (reuses shared string 392).
$ws.Cells.Item(25, 4).Value = $ws.Cells.Item(32, 4).Value()

# Resize the table (ListObject) so its definition (ref/autoFilter) covers
# the newly inserted row, matching the table growing from B2:F203 to B2:F204.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:F204"))
